# Fruta / hortaliza, semanal
# Weekly update: a new daily price record for "Clementina" / "Primera"
# (07-17-2023, serial date 45124) is inserted as row 422, pushing the
# existing rows 422:487 down to 423:488 (dimension grows from T487 to T488).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 422 - shifts rows 422:487 down to 423:488 and
# extends the used range to A1:T488 automatically.
$ws.Rows(422).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(422, 1).Value  = 5
$ws.Cells.Item(422, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(422, 3).Value  = "Maule"
$ws.Cells.Item(422, 4).Value  = 45124
$ws.Cells.Item(422, 5).Value  = 7
$ws.Cells.Item(422, 6).Value  = "Fruta"
$ws.Cells.Item(422, 7).Value  = 100102
$ws.Cells.Item(422, 8).Value  = "Cítricos"
$ws.Cells.Item(422, 9).Value  = 100102004
$ws.Cells.Item(422, 10).Value = "Mandarina"
$ws.Cells.Item(422, 11).Value = "Clementina"
$ws.Cells.Item(422, 12).Value = "Primera"
$ws.Cells.Item(422, 13).Value = 350
$ws.Cells.Item(422, 14).Value = 8000
$ws.Cells.Item(422, 15).Value = 8000
$ws.Cells.Item(422, 16).Value = 8000
$ws.Cells.Item(422, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(422, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(422, 19).Value = 800
$ws.Cells.Item(422, 20).Value = 10
